$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "61.008.86"
$ws.Range("E2").Value = "  +0.16%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.924.76"
$ws.Range("E3").Value = "  +0.22%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
Set-TextValue "D5" "591.16"
$ws.Range("E5").Value = "  +0.93%  "

# Row 6 - Solana
Set-TextValue "D6" "147.30"
$ws.Range("E6").Value = "  +1.25%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.03%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +0.36%  "

# Row 9 - Toncoin
$ws.Range("E9").Value = "  -0.24%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.53%  "

# Row 11 - Cardano
Set-TextValue "D11" "0.441"
$ws.Range("E11").Value = "  -1.40%  "

# Row 12 - ShibaInu
$ws.Range("E12").Value = "  +0.21%  "

# Row 13 - Avalanche
$ws.Range("E13").Value = "  +0.30%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +0.10%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.408.31"
$ws.Range("E15").Value = "  +0.21%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "60.943.52"
$ws.Range("E16").Value = "  +0.15%  "

# Rows 17-18 swap: WrappedEther and Polkadot swap positions
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D17" "6.71"
$ws.Range("E17").Value = "  -0.88%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.925.50"
$ws.Range("E18").Value = "  +0.20%  "

# Row 19 - BitcoinCash
Set-TextValue "D19" "432.77"
$ws.Range("E19").Value = "  +0.56%  "

# Row 20 - Chainlink
Set-TextValue "D20" "13.44"
$ws.Range("E20").Value = "  -1.45%  "

# Row 21 - Polygon
$ws.Range("E21").Value = "  -0.46%  "

# Row 22 - Uniswap
Set-TextValue "D22" "7.10"
$ws.Range("E22").Value = "  -0.49%  "

# Row 23 - Litecoin
Set-TextValue "D23" "81.39"
$ws.Range("E23").Value = "  +1.22%  "

# Row 24 - RenderToken
Set-TextValue "D24" "10.87"
$ws.Range("E24").Value = "  +0.67%  "

# Row 25 - Fetch.AI
$ws.Range("E25").Value = "  -0.21%  "

# Row 26 - InternetComputer(DFINITY)
Set-TextValue "D26" "11.92"
$ws.Range("E26").Value = "  -0.28%  "

# Row 27 - Dai
$ws.Range("E27").Value = "  +0.04%  "

# Row 28 - ImmutableX
Set-TextValue "D28" "2.27"
$ws.Range("E28").Value = "  +5.28%  "

# Row 29 - PancakeSwap
Set-TextValue "D29" "2.61"
$ws.Range("E29").Value = "  +0.04%  "

# Row 30 - NEARProtocol
Set-TextValue "D30" "7.01"
$ws.Range("E30").Value = "  -2.74%  "

# Rows 31-32 swap: EthereumClassic and Hedera swap positions
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D31" "0.110"
$ws.Range("E31").Value = "  +2.90%  "

$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D32" "26.68"
$ws.Range("E32").Value = "  +0.40%  "

# Row 33 - FirstDigitalUSD
$ws.Range("E33").Value = "  +0.04%  "

# Row 34 - PEPE
$ws.Range("E34").Value = "  -0.84%  "

# Row 35 - Mantle
$ws.Range("E35").Value = "  -0.12%  "

# Row 36 - Filecoin
$ws.Range("E36").Value = "  -0.10%  "

# Row 37 - dogwifhat
$ws.Range("E37").Value = "  -0.62%  "

# Row 38 - Stacks
$ws.Range("E38").Value = "  -1.00%  "

# Row 40 - Cosmos
$ws.Range("E40").Value = "  -1.15%  "

# Row 41 - Arweave
Set-TextValue "D41" "41.48"
$ws.Range("E41").Value = "  +1.50%  "

# Row 42 - TheGraph
Set-TextValue "D42" "0.283"
$ws.Range("E42").Value = "  -4.53%  "

# Row 43 - Bittensor
Set-TextValue "D43" "376.08"
$ws.Range("E43").Value = "  -0.55%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  -1.09%  "

# Row 45 - Maker
$ws.Range("D45").Value = "2.708.80"
$ws.Range("E45").Value = "  +0.71%  "

# Row 46 - Monero
Set-TextValue "D46" "134.16"
$ws.Range("E46").Value = "  +1.50%  "

# Row 48 - InjectiveProtocol
Set-TextValue "D48" "23.97"
$ws.Range("E48").Value = "  -3.64%  "

# Row 49 - Stellar
$ws.Range("E49").Value = "  -0.63%  "

# Row 50 - ThetaToken
Set-TextValue "D50" "2.00"
$ws.Range("E50").Value = "  -2.93%  "

# Row 51 - Cronos
$ws.Range("E51").Value = "  -0.16%  "
